$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A:A").EntireColumn.Delete()
